# Generate Report for Handoff
#
# The "cbe7af25-687c-4191-a915-7ed22c306881" entry has finished its
# handoff/handback cycle, so its row drops out of every sheet, and the
# still-open "b1e3693c-9008-434f-9246-cd16926eb1aa" entry flips from
# "Handed back: in sync with en-US" to "Ready for handoff" with refreshed
# handoff timestamps.

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = "Ready for handoff"
$ovw.Range("C2").Value = "Ready for handoff"
$ovw.Range("D2").Value = "2016-03-23 08:54:40"

# Rebuild the hyperlinks collection keeping only the row-2 link (the
# shim's Hyperlinks.Delete only operates sheet-wide, so drop everything
# and re-add what should survive, in original order, so relationship ids
# line back up with rId2, rId3, ...).
$ovw.Hyperlinks.Delete()
$ovw.Hyperlinks.Add($ovw.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cc9ee3933692a70f6a0b13275a27dbdc0851716e/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.md") | Out-Null

$ovw.Rows.Item(3).Delete()

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-23 08:54:36"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cc9ee3933692a70f6a0b13275a27dbdc0851716e/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c5bdb0691393a4c0441dbf7195855b4d41f4ea70/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/5624adc9e3139f84075fcdf6818ef4225f64f894/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.md") | Out-Null
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e56f238631cef3ac778a5448369854532c20911a/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.zh-cn.xlf") | Out-Null

$zhcn.Rows.Item(3).Delete()

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-23 08:54:40"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/cc9ee3933692a70f6a0b13275a27dbdc0851716e/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("D2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/8328874d47cfbf70fe4fb74f8221fff2c07868e7/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf") | Out-Null
$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/4e0abbaeaacb4984e498c10ca19792c609b48bd4/e2e/b1e3693c-9008-434f-9246-cd16926eb1aa.md", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.md") | Out-Null
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/24ae427b664128de15acf114b90c5414395ba7e9/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf", [Type]::Missing, [Type]::Missing, "b1e3693c-9008-434f-9246-cd16926eb1aa.263448d73583d788a29ab2ebfc86ba38fb7ef971.de-de.xlf") | Out-Null

$dede.Rows.Item(3).Delete()
